# debug status API in demo server
#
# The "RegRptFilePath" / "C:\Work\RegReport.rtf" debug row on the "stress"
# sheet is removed, rows below it shift up, and the view/selection state of
# a couple of sheets is updated (the "regression" sheet becomes the active
# tab instead of "logon").

$wb = $excel.ActiveWorkbook

# Remove the debug "RegRptFilePath" row from the "stress" sheet. This is
# row 9 (key/value pair "RegRptFilePath" / "C:\Work\RegReport.rtf"); deleting
# it shifts the following rows up and lets Excel drop the now-unused shared
# strings on save.
$stress = $wb.Worksheets.Item("stress")
$stress.Activate()
$stress.Rows(9).Delete()
$stress.Range("C17").Select()

# Make "regression" the active sheet/tab again, with a new selection.
$regression = $wb.Worksheets.Item("regression")
$regression.Activate()
$regression.Range("Q7").Select()
